# Journal de bord PreTPI - "Last step of the sprint 2"
# Fills in the next five journal rows (22-26) of the "Page 1" sheet with the
# tasks completed, their date and the time spent, and updates the current
# selection/scroll position to reflect where the author left off editing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 22 - 15/02/2021
$ws.Range("C22").Value = "Arborecense du site"
$ws.Range("D22").Value = 44242
$ws.Range("E22").Value = "120 minutes"

# Row 23 - 18/02/2021 (wraps onto two lines)
$ws.Range("C23").Value = "Transcription du template HTML en structure MVC avec redirection"
$ws.Range("D23").Value = 44245
$ws.Range("E23").Value = "120 minutes"
$ws.Rows.Item(23).RowHeight = 30

# Row 24 - 18/02/2021
$ws.Range("C24").Value = "Review des points de la docs"
$ws.Range("D24").Value = 44245
$ws.Range("E24").Value = "20 minutes"

# Row 26 - 19/02/2021 (wraps onto two lines, entered before row 25's text)
$ws.Range("C26").Value = "Regler le bug de mise en page du register"
$ws.Range("D26").Value = 44246
$ws.Range("E26").Value = "60 minutes"
$ws.Rows.Item(26).RowHeight = 30

# Row 25 - 18/02/2021 (wraps onto two lines)
$ws.Range("C25").Value = "Corriger les bugs de la base de donnée "
$ws.Range("D25").Value = 44245
$ws.Range("E25").Value = "60 minutes"
$ws.Rows.Item(25).RowHeight = 30

# Leave the selection where the author was last working.
$ws.Range("E27").Select()
